$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 corresponds to Fig "10.1126:scisignal.add0509_zero_fig3" -- delete the
# whole row (an "others mistake" duplicate entry), which shifts all following
# rows up by one.
$ws.Rows("17:17").EntireRow.Select()
$ws.Rows("17:17").Delete()

# Column A width auto-fits to its (now shorter) longest value.
$ws.Columns("A:A").EntireColumn.AutoFit()

# Leave the view/selection where the deletion happened.
$excel.ActiveWindow.ScrollRow = 12
$ws.Rows("17:17").EntireRow.Select()
